$wb = $excel.ActiveWorkbook

# Helper: set a cell to a text value while preserving its original
# (bold/centered/bordered) cell style, by writing the value with a
# leading apostrophe (forcing text) and then re-pasting the formatting
# from a neighboring header cell that already has the desired style.
function Set-HeaderText {
    param($ws, [string]$targetCellAddr, [string]$styleSourceCellAddr, [string]$text)

    $ws.Range($targetCellAddr).Value = "'" + $text
    $ws.Range($styleSourceCellAddr).Copy()
    $ws.Range($targetCellAddr).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# Sheet 1: "Potencia Acumulada - SIN (MW)"
$ws1 = $wb.Worksheets.Item(1)
Set-HeaderText $ws1 "E1" "D1" "2050"
$ws1.Rows(13).Delete()

# Sheet 2: "Geracao Periodo Medio (MWMed)"
$ws2 = $wb.Worksheets.Item(2)
Set-HeaderText $ws2 "E1" "D1" "2050"
$ws2.Rows(13).Delete()

# Sheet 3: "Atendimento a Ponta(MW)"
$ws3 = $wb.Worksheets.Item(3)
Set-HeaderText $ws3 "E1" "D1" "2050"
$ws3.Rows(13).Delete()

# Sheet 4: "Potencia Incremental - SIN(MW)" (uses decade ranges, not single years)
$ws4 = $wb.Worksheets.Item(4)
Set-HeaderText $ws4 "E1" "D1" "2041-2050"
$ws4.Rows(13).Delete()

# Sheet 5: "Emissoes Totais (MtCO2eq)" (no Total row to remove here)
$ws5 = $wb.Worksheets.Item(5)
Set-HeaderText $ws5 "E1" "D1" "2050"

# Sheet 6: "Custo Total (bilhões de R$)" - just drop the Total row
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows(4).Delete()
